$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet ---
$ws.Name = "unAutreNom"

# --- Row 1 : running totals 0..12 (A1 literal 0, then B1:M1 "=1+<prev>") ---
$ws.Range("A1").Value = 0
$ws.Range("B1").Formula = "=1+A1"
$ws.Range("C1").Formula = "=1+B1"
$ws.Range("D1").Formula = "=1+C1"
$ws.Range("E1").Formula = "=1+D1"
$ws.Range("F1").Formula = "=1+E1"
$ws.Range("G1").Formula = "=1+F1"
$ws.Range("H1").Formula = "=1+G1"
$ws.Range("I1").Formula = "=1+H1"
$ws.Range("J1").Formula = "=1+I1"
$ws.Range("K1").Formula = "=1+J1"
$ws.Range("L1").Formula = "=1+K1"
$ws.Range("M1").Formula = "=1+L1"

# --- Column A : row counter, rows 2..9 ---
$ws.Range("A2").Formula = "=A1+1"
$ws.Range("A3").Formula = "=A2+1"
$ws.Range("A4").Formula = "=A3+1"
$ws.Range("A5").Formula = "=A4+1"
$ws.Range("A6").Formula = "=A5+1"
$ws.Range("A7").Formula = "=A6+1"
$ws.Range("A8").Formula = "=A7+1"
$ws.Range("A9").Formula = "=A8+1"

# --- Row 2 : header labels (B2:K2) ---
$ws.Range("B2").Value = "id"
$ws.Range("C2").Value = "NOM"
$ws.Range("D2").Value = "JOUEUR"
$ws.Range("E2").Value = "PROFIL"
$ws.Range("F2").Value = "NIVEAU"
$ws.Range("G2").Value = "RACE"
$ws.Range("H2").Value = "SEXE"
$ws.Range("I2").Value = "AGE"
$ws.Range("J2").Value = "TAILLE"
$ws.Range("K2").Value = "POIDS"

# --- Formatting -----------------------------------------------------------
# Build each distinct style on a single anchor cell first (single-cell
# writes settle directly on the final style with no stray intermediates),
# then fan the finished style out with Copy / PasteSpecial(xlPasteFormats)
# so the rest of the range picks up the very same style index instead of
# re-deriving (and leaving behind unused) intermediate styles.

# Style A: centered (horizontal + vertical), regular weight -> anchor A1
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

# Style B: centered + bold -> anchor B2
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4108

# Fan style A out to the rest of column A and to row 1 (B1:M1)
$ws.Range("A1").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)
$ws.Range("B1:M1").PasteSpecial(-4122)

# Fan style B out to the rest of the header row (C2:K2)
$ws.Range("B2").Copy()
$ws.Range("C2:K2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column widths ---
# (ColumnWidth is internally quantized to 1/6-character steps by this
# engine, so these values are chosen to land on the closest reachable
# width to the original 4.7109375 / 10.42578125 character widths.)
$ws.Columns("A").ColumnWidth = 3.8333333333333
$ws.Range("B1:K1").EntireColumn.ColumnWidth = 9.6666666666667

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$null = $ws.Range("B3").Select()
